$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New note appended to Sheet1 (written first so this shared string lands
# at index 33, ahead of everything that gets typed onto the new sheet).
$ws1.Range("A12").Value = "TimeLog for Assignment 2 to on List 2"

# New sheet "Лист2", inserted right after "Лист1" -> becomes sheetId 2 /
# the active tab (activeTab=1).
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Лист2"

# Header row: reuse Sheet1's header formatting (yellow fill style) instead
# of re-creating an equivalent style, so styles.xml stays untouched.
$ws1.Range("A1:C1").Copy() | Out-Null
$ws2.Range("A1:C1").PasteSpecial(-4122) | Out-Null

$ws2.Range("A1").Value = "Activity"
$ws2.Range("B1").Value = "Estimated Time"
$ws2.Range("C1").Value = "Real Time"

$ws2.Range("A3").Value = "Study the book chapters 4-5"
$ws2.Range("B3").Value = "2 hours "
$ws2.Range("C3").Value = "2.5 hours"

$ws2.Range("A5").Value = "Study the book chapters 6-7"
$ws2.Range("B5").Value = "2 hours "
$ws2.Range("C5").Value = "1.5 hours"

$ws2.Range("A7").Value = "Create Use Case Model"
$ws2.Range("B7").Value = "4 hours"
$ws2.Range("C7").Value = "8 hours"

$ws2.Range("A9").Value = "Draw State Machine Diagram"
$ws2.Range("B9").Value = "2 hours"
$ws2.Range("C9").Value = "3 hours"

$ws2.Range("A11").Value = "Draw Class Diagram"
$ws2.Range("B11").Value = "2 hours"
$ws2.Range("C11").Value = "1 hour"

$ws2.Range("A13").Value = "Implement The Game / Refactor"
$ws2.Range("B13").Value = "6 hours"
$ws2.Range("C13").Value = "8 hours"

# Hide the (unused) D:E columns on the new sheet, like the authored file.
$ws2.Range("D1:E1").EntireColumn.Hidden = $true

# Column width tweaks. ColumnWidth is funnelled through Excel's
# integer-pixel rounding on its way to the stored OOXML "width", so the
# raw values below are chosen so the rounded pixel width lands as close
# as that quantisation allows to the authored figures (33.62/12.5 on
# Sheet1, 27.23/22.23/14.31 on Sheet2 - 12.5 lands exactly, the others
# land on the nearest achievable pixel boundary).
$ws1.Columns.Item(1).ColumnWidth = 32.833333333333336
$ws1.Columns.Item(4).ColumnWidth = 11.666666666666666

$ws2.Columns.Item(1).ColumnWidth = 26.333333333333332
$ws2.Columns.Item(2).ColumnWidth = 21.333333333333332
$ws2.Columns.Item(3).ColumnWidth = 13.5

# Match the authored page setup for the new sheet (margins + header/footer
# text mirror Sheet1's).
$ps2 = $ws2.PageSetup
$ps2.LeftMargin = 56.7
$ps2.RightMargin = 56.7
$ps2.TopMargin = 75.8
$ps2.BottomMargin = 75.8
$ps2.HeaderMargin = 56.7
$ps2.FooterMargin = 56.7
$ps2.CenterHeader = "&`"Times New Roman,Обычный`"&12&A"
$ps2.CenterFooter = "&`"Times New Roman,Обычный`"&12Страница &P"

# Selections: Sheet1 is no longer the active tab, cursor parked at D16;
# Sheet2 becomes the active tab with the cursor at A9.
$ws1.Activate() | Out-Null
$ws1.Range("D16").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A9").Select() | Out-Null
